$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from an existing header cell (H1) so new header cells
# I1 and J1 share the same bold/centered/bordered style (s="1").
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iVals = @(8,8,7,10,9,7,7,8,7,7,8,8,9,8,9,9,9,9,8,7,8,9,9,8,8,8,10,8,9,9,9,9,9,9,9,9,9,9,8,9,9,9,8,8,9,8,8,7,9,9,9,9,9,8,9,9,9,8,8,8,9,9,9,9,8,9,5,8)
$jVals = @(8,8,7,10,9,8,7,8,7,7,8,8,9,8,9,9,9,9,9,7,8,9,9,9,8,8,10,8,9,9,9,9,9,9,10,9,9,9,9,9,10,9,8,8,9,8,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,5,8)

for ($i = 0; $i -lt 68; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$i]
    $ws.Cells.Item($r, 10).Value = $jVals[$i]
}
